$d = $word.ActiveDocument

# Locate the end of the existing "about me" paragraph and split it into:
#   1) the original sentence (unchanged)
#   2) a new blank paragraph
#   3) a new paragraph: "I love eating and I also love testcode cap "
#   4) a trailing blank paragraph (keeps the _GoBack bookmark that sat at
#      the end of the original paragraph)
$find = $d.Content.Find
$replaced = $find.Execute(
    "technology;",                                                      `
    $true,                                                               `
    $false,                                                              `
    $false,                                                              `
    $false,                                                              `
    $false,                                                              `
    $true,                                                               `
    1,                                                                   `
    $false,                                                              `
    "technology;^p^pI love eating and I also love testcode cap ^p",      `
    2                                                                    `
)

if (-not $replaced) {
    throw "Could not find the anchor text to split the paragraph."
}
